$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2.96768027051047
$ws.Range("D2").Value = 9.728771743612775
$ws.Range("E2").Value = 13.20954494546152
$ws.Range("F2").Value = 27.76872353642566
$ws.Range("G2").Value = 27.98682674982453
$ws.Range("H2").Value = 13.48843247488505
$ws.Range("I2").Value = 19.31845659362894
$ws.Range("J2").Value = 9.46981976747378
$ws.Range("M2").Value = 21.72846864070457
$ws.Range("N2").Value = 17.47360446468254
$ws.Range("O2").Value = 20.54381159118376
$ws.Range("C3").Value = 2.882449125003153
$ws.Range("D3").Value = 9.768708863714499
$ws.Range("E3").Value = 13.29089653459492
$ws.Range("F3").Value = 27.66892435197116
$ws.Range("G3").Value = 27.55023661399635
$ws.Range("H3").Value = 13.48733439115702
$ws.Range("I3").Value = 19.22560426770237
$ws.Range("J3").Value = 9.515145135360529
$ws.Range("M3").Value = 21.05999307995793
$ws.Range("N3").Value = 17.17474433459105
$ws.Range("O3").Value = 20.46115373467746
$ws.Range("C4").Value = 2.828186969592786
$ws.Range("D4").Value = 9.794701107102281
$ws.Range("E4").Value = 13.34335068796662
$ws.Range("F4").Value = 27.6158772931579
$ws.Range("G4").Value = 27.28896514734069
$ws.Range("H4").Value = 13.48958019405705
$ws.Range("I4").Value = 19.17384359585879
$ws.Range("J4").Value = 9.544312427348959
$ws.Range("M4").Value = 20.63826577681662
$ws.Range("N4").Value = 16.99048499217523
$ws.Range("O4").Value = 20.41598820084439
$ws.Range("C5").Value = 2.805607131517577
$ws.Range("D5").Value = 9.805663278172352
$ws.Range("E5").Value = 13.36535711874237
$ws.Range("F5").Value = 27.59634497980137
$ws.Range("G5").Value = 27.18435796623889
$ws.Range("H5").Value = 13.49122995981464
$ws.Range("I5").Value = 19.15408824526617
$ws.Range("J5").Value = 9.556535555646958
$ws.Range("M5").Value = 20.46379525394675
$ws.Range("N5").Value = 16.91530027709114
$ws.Range("O5").Value = 20.39900144105933
$ws.Range("C6").Value = 2.801830075210763
$ws.Range("D6").Value = 9.807505901875063
$ws.Range("E6").Value = 13.3690494152917
$ws.Range("F6").Value = 27.59322795943521
$ws.Range("G6").Value = 27.16710497064128
$ws.Range("H6").Value = 13.49154827063247
$ws.Range("I6").Value = 19.15088912824666
$ws.Range("J6").Value = 9.558585591586841
$ws.Range("M6").Value = 20.43467362032846
$ws.Range("N6").Value = 16.9028127771975
$ws.Range("O6").Value = 20.39626684196474
$ws.Range("C7").Value = 2.82788431861117
$ws.Range("D7").Value = 9.794847447757862
$ws.Range("E7").Value = 13.34364491802651
$ws.Range("F7").Value = 27.61560541434108
$ws.Range("G7").Value = 27.28754663431188
$ws.Range("H7").Value = 13.48959946898073
$ws.Range("I7").Value = 19.17357173158505
$ws.Range("J7").Value = 9.544475906107666
$ws.Range("M7").Value = 20.63592307827573
$ws.Range("N7").Value = 16.98947129541588
$ws.Range("O7").Value = 20.41575335115297
$ws.Range("C8").Value = 2.938703261937857
$ws.Range("D8").Value = 9.742237003196147
$ws.Range("E8").Value = 13.2370760632535
$ws.Range("F8").Value = 27.73261468403218
$ws.Range("G8").Value = 27.83497330782886
$ws.Range("H8").Value = 13.48744810773176
$ws.Range("I8").Value = 19.28536216427558
$ws.Range("J8").Value = 9.485170957722261
$ws.Range("M8").Value = 21.50044776447563
$ws.Range("N8").Value = 17.37077004203107
$ws.Range("O8").Value = 20.51416048319687
$ws.Range("C9").Value = 3.140007479825398
$ws.Range("D9").Value = 9.650723288290584
$ws.Range("E9").Value = 13.04790217089713
$ws.Range("F9").Value = 28.02656401193218
$ws.Range("G9").Value = 28.95548819400575
$ws.Range("H9").Value = 13.50636908956668
$ws.Range("I9").Value = 19.54539827266024
$ws.Range("J9").Value = 9.37944173035563
$ws.Range("M9").Value = 23.09740890401267
$ws.Range("N9").Value = 18.10850323604079
$ws.Range("O9").Value = 20.75076860911644
$ws.Range("C10").Value = 3.277316211266564
$ws.Range("D10").Value = 9.590573575779514
$ws.Range("E10").Value = 12.92090512042861
$ws.Range("F10").Value = 28.28057446514243
$ws.Range("G10").Value = 29.79804656530331
$ws.Range("H10").Value = 13.53431480321853
$ws.Range("I10").Value = 19.76005955162619
$ws.Range("J10").Value = 9.30814268993776
$ws.Range("M10").Value = 24.20016031821839
$ws.Range("N10").Value = 18.63929182088354
$ws.Range("O10").Value = 20.95014269368189
$ws.Range("C11").Value = 3.337335263754881
$ws.Range("D11").Value = 9.564744462909156
$ws.Range("E11").Value = 12.865716063064
$ws.Range("F11").Value = 28.40406427914127
$ws.Range("G11").Value = 30.18357872464424
$ws.Range("H11").Value = 13.5500557788445
$ws.Range("I11").Value = 19.86253427555567
$ws.Range("J11").Value = 9.277079047430883
$ws.Range("M11").Value = 24.68463481784001
$ws.Range("N11").Value = 18.87730401731017
$ws.Range("O11").Value = 21.04612480188285
$ws.Range("C12").Value = 3.359700987927871
$ws.Range("D12").Value = 9.555183902348203
$ws.Range("E12").Value = 12.8451874925342
$ws.Range("F12").Value = 28.4519368812406
$ws.Range("G12").Value = 30.32973339364705
$ws.Range("H12").Value = 13.55644934066202
$ws.Range("I12").Value = 19.90200425241133
$ws.Range("J12").Value = 9.265512172211666
$ws.Range("M12").Value = 24.86548866751528
$ws.Range("N12").Value = 18.96685761791887
$ws.Range("O12").Value = 21.08320610077721
$ws.Range("C13").Value = 3.354900397232762
$ws.Range("D13").Value = 9.557233141143705
$ws.Range("E13").Value = 12.84959223331523
$ws.Range("F13").Value = 28.44157783905959
$ws.Range("G13").Value = 30.29825168208575
$ws.Range("H13").Value = 13.55505317119706
$ws.Range("I13").Value = 19.89347458081606
$ws.Range("J13").Value = 9.267994590775153
$ws.Range("M13").Value = 24.82665647825059
$ws.Range("N13").Value = 18.94759752804431
$ws.Range("O13").Value = 21.07518770796688
$ws.Range("C14").Value = 3.33918262192198
$ws.Range("D14").Value = 9.563953494054701
$ws.Range("E14").Value = 12.86401974914271
$ws.Range("F14").Value = 28.40798073823486
$ws.Range("G14").Value = 30.19560061715946
$ws.Range("H14").Value = 13.5505731240068
$ws.Range("I14").Value = 19.86576833542357
$ws.Range("J14").Value = 9.276123505281756
$ws.Range("M14").Value = 24.69956673242186
$ws.Range("N14").Value = 18.88468364703407
$ws.Range("O14").Value = 21.04916092222202
$ws.Range("C15").Value = 3.329507562114668
$ws.Range("D15").Value = 9.568098598815494
$ws.Range("E15").Value = 12.87290521642647
$ws.Range("F15").Value = 28.38754507077883
$ws.Range("G15").Value = 30.1327402683363
$ws.Range("H15").Value = 13.54788524178256
$ws.Range("I15").Value = 19.84888319089672
$ws.Range("J15").Value = 9.281128232963091
$ws.Range("M15").Value = 24.62137740235497
$ws.Range("N15").Value = 18.84606967927996
$ws.Range("O15").Value = 21.03331370638283
$ws.Range("C16").Value = 3.273343704184127
$ws.Range("D16").Value = 9.592292412040312
$ws.Range("E16").Value = 12.92456373303235
$ws.Range("F16").Value = 28.272661193308
$ws.Range("G16").Value = 29.77288247971565
$ws.Range("H16").Value = 13.53334681250541
$ws.Range("I16").Value = 19.75345724188922
$ws.Range("J16").Value = 9.310200280543418
$ws.Range("M16").Value = 24.16814117231606
$ws.Range("N16").Value = 18.62366080583425
$ws.Range("O16").Value = 20.94397428686265
$ws.Range("C17").Value = 3.238255091955496
$ws.Range("D17").Value = 9.60752715407963
$ws.Range("E17").Value = 12.95691526256864
$ws.Range("F17").Value = 28.20419507744638
$ws.Range("G17").Value = 29.55258102385211
$ws.Range("H17").Value = 13.52520182045004
$ws.Range("I17").Value = 19.69613214323925
$ws.Range("J17").Value = 9.328385492254446
$ws.Range("M17").Value = 23.88559487677704
$ws.Range("N17").Value = 18.4862786156291
$ws.Range("O17").Value = 20.89050386945996
$ws.Range("C18").Value = 3.217843635305888
$ws.Range("D18").Value = 9.616434101794615
$ws.Range("E18").Value = 12.97576614771843
$ws.Range("F18").Value = 28.16556451805383
$ws.Range("G18").Value = 29.42609045396533
$ws.Range("H18").Value = 13.52080226485005
$ws.Range("I18").Value = 19.66361615323029
$ws.Range("J18").Value = 9.33897420705685
$ws.Range("M18").Value = 23.72147250240288
$ws.Range("N18").Value = 18.40694021676022
$ws.Range("O18").Value = 20.86024823177984
$ws.Range("C19").Value = 3.210893614646826
$ws.Range("D19").Value = 9.619474633540637
$ws.Range("E19").Value = 12.98219052532207
$ws.Range("F19").Value = 28.15261450978504
$ws.Range("G19").Value = 29.38330586543428
$ws.Range("H19").Value = 13.51936171356655
$ws.Range("I19").Value = 19.65268594889069
$ws.Range("J19").Value = 9.342581556366019
$ws.Range("M19").Value = 23.66563152183276
$ws.Range("N19").Value = 18.3800252069867
$ws.Range("O19").Value = 20.85009066056266
$ws.Range("C20").Value = 3.242014160528364
$ws.Range("D20").Value = 9.605890453407975
$ws.Range("E20").Value = 12.95344623042698
$ws.Range("F20").Value = 28.21140608262968
$ws.Range("G20").Value = 29.57601080002762
$ws.Range("H20").Value = 13.52603936665736
$ws.Range("I20").Value = 19.70218751835554
$ws.Range("J20").Value = 9.326436294125399
$ws.Range("M20").Value = 23.91583999952654
$ws.Range("N20").Value = 18.50093688153605
$ws.Range("O20").Value = 20.89614440570923
$ws.Range("C21").Value = 3.343809223606537
$ws.Range("D21").Value = 9.561973583775615
$ws.Range("E21").Value = 12.85977199556585
$ws.Range("F21").Value = 28.41781917046113
$ws.Range("G21").Value = 30.22574858801426
$ws.Range("H21").Value = 13.55187729996599
$ws.Range("I21").Value = 19.87388851815739
$ws.Range("J21").Value = 9.273730526523639
$ws.Range("M21").Value = 24.7369678149063
$ws.Range("N21").Value = 18.90317924715989
$ws.Range("O21").Value = 21.05678588116871
$ws.Range("C22").Value = 3.408222875395778
$ws.Range("D22").Value = 9.534555775606302
$ws.Range("E22").Value = 12.8007085590602
$ws.Range("F22").Value = 28.55917378668322
$ws.Range("G22").Value = 30.6512628581675
$ws.Range("H22").Value = 13.5712850036457
$ws.Range("I22").Value = 19.98996708344594
$ws.Range("J22").Value = 9.240427780135088
$ws.Range("M22").Value = 25.25837893331994
$ws.Range("N22").Value = 19.16267340801188
$ws.Range("O22").Value = 21.16604560298678
$ws.Range("C23").Value = 3.374040920346357
$ws.Range("D23").Value = 9.549071688611047
$ws.Range("E23").Value = 12.8320346971692
$ws.Range("F23").Value = 28.48315097035919
$ws.Range("G23").Value = 30.4241300006765
$ws.Range("H23").Value = 13.56069705952985
$ws.Range("I23").Value = 19.927670212561
$ws.Range("J23").Value = 9.25809772812233
$ws.Range("M23").Value = 24.98152804609281
$ws.Range("N23").Value = 19.02451269730234
$ws.Range("O23").Value = 21.10734949212162
$ws.Range("C24").Value = 3.240315427821687
$ws.Range("D24").Value = 9.606629943853758
$ws.Range("E24").Value = 12.95501379616165
$ws.Range("F24").Value = 28.20814370608652
$ws.Range("G24").Value = 29.56541767597987
$ws.Range("H24").Value = 13.52565982991406
$ws.Range("I24").Value = 19.69944850745936
$ws.Range("J24").Value = 9.327317109754969
$ws.Range("M24").Value = 23.90217141032969
$ws.Range("N24").Value = 18.49431098076904
$ws.Range("O24").Value = 20.89359280590667
$ws.Range("C25").Value = 3.087354626509685
$ws.Range("D25").Value = 9.674234520509771
$ws.Range("E25").Value = 13.09696719916833
$ws.Range("F25").Value = 27.94025869001536
$ws.Range("G25").Value = 28.64829234691243
$ws.Range("H25").Value = 13.49877792428289
$ws.Range("I25").Value = 19.47080579613275
$ws.Range("J25").Value = 9.40691939455292
$ws.Range("M25").Value = 22.67708084714791
$ws.Range("N25").Value = 17.91053905459153
$ws.Range("O25").Value = 20.68219094432716